$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($rowA, $rowB)

    $addrBC_A = "B" + $rowA + ":C" + $rowA
    $addrBC_B = "B" + $rowB + ":C" + $rowB
    $addrE_AD_A = "E" + $rowA + ":AD" + $rowA
    $addrE_AD_B = "E" + $rowB + ":AD" + $rowB

    $srcBC_A = $ws.Range($addrBC_A).Value()
    $srcBC_B = $ws.Range($addrBC_B).Value()
    $srcE_AD_A = $ws.Range($addrE_AD_A).Value()
    $srcE_AD_B = $ws.Range($addrE_AD_B).Value()

    $ws.Range($addrBC_A).Value = $srcBC_B
    $ws.Range($addrBC_B).Value = $srcBC_A
    $ws.Range($addrE_AD_A).Value = $srcE_AD_B
    $ws.Range($addrE_AD_B).Value = $srcE_AD_A
}

Swap-Rows 92 93
Swap-Rows 94 95
Swap-Rows 98 99
Swap-Rows 126 127
Swap-Rows 276 277
Swap-Rows 282 283
Swap-Rows 298 299
Swap-Rows 303 304
Swap-Rows 310 312
Swap-Rows 318 319
Swap-Rows 322 323
